$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("B3").Value = "Grand Plaza Apartments"
$ws.Range("A1:B3").Select()
